$d = $word.ActiveDocument

# 1. Update the numeric result of the example from 31.972.500. to 23.962.500.
$rng = $d.Content
$rng.Find.Execute("31.972.500.", $false, $false, $false, $false, $false, $true, 1, $false, "23.962.500.", 2)

# 2. Locate the replaced text so we know its absolute character offsets.
$rng2 = $d.Content
$rng2.Find.Execute("23.962.500.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $rng2.Start

# 3. Split "23.962.500." into separate runs: "23" "." "962" "." "500" "."
#    by inserting and immediately deleting temporary bookmarks at each
#    internal boundary (Word splits runs at bookmark insertion points and
#    leaves the split in place even after the bookmark itself is removed).
$splitOffsets = @(2, 3, 6, 10)
$i = 0
foreach ($off in $splitOffsets) {
    $i = $i + 1
    $markName = "tempSplitMark$i"
    $splitRng = $d.Range($base + $off, $base + $off)
    $bm = $d.Bookmarks.Add($markName, $splitRng)
    $d.Bookmarks($markName).Delete()
}

# 4. Place the _GoBack bookmark between "962." and "500" (offset 7), which
#    also removes it from its old location at the top of the document.
$bmRng = $d.Range($base + 7, $base + 7)
$d.Bookmarks.Add("_GoBack", $bmRng)
